$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 13, shifting the existing rows 13-15 down to 14-16.
$ws.Rows("13:13").Insert()

# Populate the new row 13 with the new weekly record. Columns A, B, C, E, F, G, H, I, R
# mirror the fixed "Cebollín" attributes shared by every row in this sheet.
$ws.Range("A13").Value = 7
$ws.Range("B13").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C13").Value = "Ñuble"
$ws.Range("D13").Value = 44762
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = 100112037
$ws.Range("G13").Value = "Cebollín"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 8000
$ws.Range("L13").Value = 8000
$ws.Range("M13").Value = 8000
$ws.Range("N13").Value = "$/docena de atados"
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 2667
$ws.Range("Q13").Value = 3
$ws.Range("R13").Value = "Hortaliza"
